$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "76-73=3"
$t.Cell(1, 2).Range.Text = "56-54=2"
$t.Cell(1, 3).Range.Text = "30+38=68"
$t.Cell(1, 4).Range.Text = "57-12=45"
$t.Cell(1, 5).Range.Text = "51+14=65"

$t.Cell(2, 1).Range.Text = "3+23=26"
$t.Cell(2, 2).Range.Text = "38+11=49"
$t.Cell(2, 3).Range.Text = "9+3=12"
$t.Cell(2, 4).Range.Text = "87-16=71"
$t.Cell(2, 5).Range.Text = "65+29=94"

$t.Cell(3, 1).Range.Text = "95-34=61"
$t.Cell(3, 2).Range.Text = "16+15=31"
$t.Cell(3, 3).Range.Text = "93-84=9"
$t.Cell(3, 4).Range.Text = "62-58=4"
$t.Cell(3, 5).Range.Text = "5+47=52"

$t.Cell(4, 1).Range.Text = "34-17=17"
$t.Cell(4, 2).Range.Text = "53-40=13"
$t.Cell(4, 3).Range.Text = "34+16=50"
$t.Cell(4, 4).Range.Text = "78-31=47"
$t.Cell(4, 5).Range.Text = "4+87=91"

$t.Cell(5, 1).Range.Text = "60-16=44"
$t.Cell(5, 2).Range.Text = "39+13=52"
$t.Cell(5, 3).Range.Text = "49-21=28"
$t.Cell(5, 4).Range.Text = "21+26=47"
$t.Cell(5, 5).Range.Text = "67-0=67"

$t.Cell(6, 1).Range.Text = "58-56=2"
$t.Cell(6, 2).Range.Text = "73-11=62"
$t.Cell(6, 3).Range.Text = "76-5=71"
$t.Cell(6, 4).Range.Text = "81-80=1"
$t.Cell(6, 5).Range.Text = "97-4=93"

$t.Cell(7, 1).Range.Text = "89-77=12"
$t.Cell(7, 2).Range.Text = "47-20=27"
$t.Cell(7, 3).Range.Text = "61-30=31"
$t.Cell(7, 4).Range.Text = "25-11=14"
$t.Cell(7, 5).Range.Text = "61-24=37"

$t.Cell(8, 1).Range.Text = "60-12=48"
$t.Cell(8, 2).Range.Text = "70+4=74"
$t.Cell(8, 3).Range.Text = "76+18=94"
$t.Cell(8, 4).Range.Text = "56+14=70"
$t.Cell(8, 5).Range.Text = "68-22=46"

$t.Cell(9, 1).Range.Text = "8+55=63"
$t.Cell(9, 2).Range.Text = "77-15=62"
$t.Cell(9, 3).Range.Text = "60+36=96"
$t.Cell(9, 4).Range.Text = "61-13=48"
$t.Cell(9, 5).Range.Text = "82-74=8"

$t.Cell(10, 1).Range.Text = "10+8=18"
$t.Cell(10, 2).Range.Text = "31-12=19"
$t.Cell(10, 3).Range.Text = "66-43=23"
$t.Cell(10, 4).Range.Text = "63-30=33"
$t.Cell(10, 5).Range.Text = "14+67=81"

$t.Cell(11, 1).Range.Text = "49-30=19"
$t.Cell(11, 2).Range.Text = "28+66=94"
$t.Cell(11, 3).Range.Text = "51-16=35"
$t.Cell(11, 4).Range.Text = "97-95=2"
$t.Cell(11, 5).Range.Text = "89-71=18"

$t.Cell(12, 1).Range.Text = "64-61=3"
$t.Cell(12, 2).Range.Text = "7+74=81"
$t.Cell(12, 3).Range.Text = "90-6=84"
$t.Cell(12, 4).Range.Text = "76-9=67"
$t.Cell(12, 5).Range.Text = "76+20=96"

$t.Cell(13, 1).Range.Text = "10-8=2"
$t.Cell(13, 2).Range.Text = "35-31=4"
$t.Cell(13, 3).Range.Text = "17+8=25"
$t.Cell(13, 4).Range.Text = "48-40=8"
$t.Cell(13, 5).Range.Text = "87-83=4"

$t.Cell(14, 1).Range.Text = "6+62=68"
$t.Cell(14, 2).Range.Text = "89-37=52"
$t.Cell(14, 3).Range.Text = "3+19=22"
$t.Cell(14, 4).Range.Text = "71-10=61"
$t.Cell(14, 5).Range.Text = "53-15=38"

$t.Cell(15, 1).Range.Text = "0-0=0"
$t.Cell(15, 2).Range.Text = "45+17=62"
$t.Cell(15, 3).Range.Text = "80+18=98"
$t.Cell(15, 4).Range.Text = "0+10=10"
$t.Cell(15, 5).Range.Text = "96-69=27"

$t.Cell(16, 1).Range.Text = "76-7=69"
$t.Cell(16, 2).Range.Text = "19+75=94"
$t.Cell(16, 3).Range.Text = "10+62=72"
$t.Cell(16, 4).Range.Text = "66-15=51"
$t.Cell(16, 5).Range.Text = "98-45=53"

$t.Cell(17, 1).Range.Text = "72+24=96"
$t.Cell(17, 2).Range.Text = "38+1=39"
$t.Cell(17, 3).Range.Text = "48+19=67"
$t.Cell(17, 4).Range.Text = "43-31=12"
$t.Cell(17, 5).Range.Text = "44+38=82"

$t.Cell(18, 1).Range.Text = "83+16=99"
$t.Cell(18, 2).Range.Text = "60-20=40"
$t.Cell(18, 3).Range.Text = "20+21=41"
$t.Cell(18, 4).Range.Text = "13+51=64"
$t.Cell(18, 5).Range.Text = "9+7=16"

$t.Cell(19, 1).Range.Text = "86-4=82"
$t.Cell(19, 2).Range.Text = "73-37=36"
$t.Cell(19, 3).Range.Text = "0+79=79"
$t.Cell(19, 4).Range.Text = "96-27=69"
$t.Cell(19, 5).Range.Text = "60-47=13"

$t.Cell(20, 1).Range.Text = "38+61=99"
$t.Cell(20, 2).Range.Text = "95-31=64"
$t.Cell(20, 3).Range.Text = "44-32=12"
$t.Cell(20, 4).Range.Text = "75-18=57"
$t.Cell(20, 5).Range.Text = "7+75=82"
